$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: municipio-nombre and aragon dimensions are now curated as refArea
$ws.Range("C2").Value = "sdmx-dimension:refArea"
$ws.Range("E2").Value = "sdmx-dimension:refArea"

# Row 3: municipio-nombre is no longer a measure ("medida") but a dimension ("dim")
$ws.Range("C3").Value = "dim"

# Row 4: new URI types for the curated dimensions
$ws.Range("C4").Value = "URI-Municipio"
$ws.Range("E4").Value = "URI-Comunidad"

# Row 5: the aragon-specific mapping file is no longer used
$ws.Range("E5").Value = $null
